$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws "D2" "333.60"
Set-TextValue $ws "E2" "0.86%"

# Row 3
Set-TextValue $ws "D3" "39.50"
Set-TextValue $ws "E3" "-2.13%"

# Row 4
Set-TextValue $ws "D4" "5.759"
Set-TextValue $ws "E4" "2.72%"

# Row 5
Set-TextValue $ws "D5" "0.08064"

# Row 6
Set-TextValue $ws "D6" "2.007"
Set-TextValue $ws "E6" "3.92%"

# Row 7
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue $ws "D7" "8.665"
Set-TextValue $ws "E7" "-0.18%"

# Row 8
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws "D8" "4.493"
Set-TextValue $ws "E8" "-1.17%"

# Row 9
Set-TextValue $ws "D9" "2.963"
Set-TextValue $ws "E9" "-0.48%"

# Row 10
Set-TextValue $ws "D10" "0.9226"
Set-TextValue $ws "E10" "-2.76%"

# Row 11
Set-TextValue $ws "D11" "0.1278"
Set-TextValue $ws "E11" "1.24%"

# Row 12
Set-TextValue $ws "D12" "0.1963"
Set-TextValue $ws "E12" "-0.89%"

# Row 13
Set-TextValue $ws "D13" "8.744"
Set-TextValue $ws "E13" "18.52%"

# Row 14
Set-TextValue $ws "D14" "0.09255"
Set-TextValue $ws "E14" "0.29%"

# Row 15
Set-TextValue $ws "D15" "0.03549"
Set-TextValue $ws "E15" "0.29%"

# Row 16
Set-TextValue $ws "D16" "0.1051"
Set-TextValue $ws "E16" "9.65%"

# Row 17
Set-TextValue $ws "D17" "0.001307"
Set-TextValue $ws "E17" "-1.68%"

# Row 18
Set-TextValue $ws "D18" "0.006381"
Set-TextValue $ws "E18" "1.73%"

# Row 19
Set-TextValue $ws "D19" "3.366"
Set-TextValue $ws "E19" "0.02%"

# Row 20
Set-TextValue $ws "D20" "0.3485"
Set-TextValue $ws "E20" "-0.80%"

# Row 21
Set-TextValue $ws "D21" "0.1360"
Set-TextValue $ws "E21" "2.04%"

# Row 22
Set-TextValue $ws "D22" "0.2717"
Set-TextValue $ws "E22" "10.83%"

# Row 23
Set-TextValue $ws "D23" "0.04418"
Set-TextValue $ws "E23" "-0.33%"

# Row 24
Set-TextValue $ws "D24" "0.001259"
Set-TextValue $ws "E24" "2.94%"

# Row 25
Set-TextValue $ws "D25" "0.004532"
Set-TextValue $ws "E25" "4.90%"

# Row 26
Set-TextValue $ws "D26" "0.0001148"
Set-TextValue $ws "E26" "-4.29%"

# Row 39
Set-TextValue $ws "D39" "0.02535"
Set-TextValue $ws "E39" "0.41%"

# Row 40
Set-TextValue $ws "D40" "0.05488"
Set-TextValue $ws "E40" "4.58%"

# Row 41
Set-TextValue $ws "D41" "0.007449"
Set-TextValue $ws "E41" "-5.70%"

# Row 42
Set-TextValue $ws "D42" "0.009896"
Set-TextValue $ws "E42" "16.58%"

# Row 43
Set-TextValue $ws "D43" "0.1413"
Set-TextValue $ws "E43" "-1.32%"

# Row 44
Set-TextValue $ws "D44" "0.002105"
Set-TextValue $ws "E44" "-1.65%"

# Row 45
Set-TextValue $ws "D45" "0.01128"
Set-TextValue $ws "E45" "8.69%"

# Row 46
Set-TextValue $ws "D46" "0.00007258"
Set-TextValue $ws "E46" "9.87%"

# Row 47
Set-TextValue $ws "D47" "0.00000000749"
Set-TextValue $ws "E47" "-0.13%"

# Row 48
Set-TextValue $ws "D48" "0.003025"
Set-TextValue $ws "E48" "4.09%"

# Row 49
Set-TextValue $ws "D49" "0.002276"
Set-TextValue $ws "E49" "-5.16%"

# Row 50
Set-TextValue $ws "D50" "0.00002097"
Set-TextValue $ws "E50" "-0.13%"

# Row 51
Set-TextValue $ws "D51" "0.0001998"
Set-TextValue $ws "E51" "-0.13%"
